$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 0.68

# Row 6
$ws.Range("D6").Value = 1.55
$ws.Range("F6").Value = 1.17
$ws.Range("G6").Value = 1.04

# Row 7
$ws.Range("C7").Value = 2.11
$ws.Range("F7").Value = 1.46
